$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared-string rich text cells) ---
$ws.Range("A8").Value = "Volume 30   Number  18"
$ws.Range("C9").Value = "Report Covering the Week  5/1/2023  Through  5/7/2023"

# --- Weekly crime statistics table updates (rows 14-30) ---
# Row 14
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 5
$ws.Range("G14").Value = 3
$ws.Range("H14").Value = 66.666666666666
$ws.Range("I14").Value = 28
$ws.Range("J14").Value = 22
$ws.Range("K14").Value = 27.272727272727
$ws.Range("L14").Value = 33.333333333333
$ws.Range("M14").Value = 12
$ws.Range("N14").Value = -64.102564102564

# Row 15
$ws.Range("C15").Value = 2
$ws.Range("D15").Value = 5
$ws.Range("E15").Value = -60
$ws.Range("F15").Value = 9
$ws.Range("G15").Value = 12
$ws.Range("H15").Value = -25
$ws.Range("I15").Value = 71
$ws.Range("J15").Value = 75
$ws.Range("K15").Value = -5.333333333333
$ws.Range("L15").Value = 9.230769230769
$ws.Range("M15").Value = 39.215686274509
$ws.Range("N15").Value = -60.989010989011

# Row 16
$ws.Range("C16").Value = 41
$ws.Range("D16").Value = 31
$ws.Range("E16").Value = 32.258064516129
$ws.Range("F16").Value = 142
$ws.Range("G16").Value = 118
$ws.Range("H16").Value = 20.338983050847
$ws.Range("I16").Value = 611
$ws.Range("J16").Value = 652
$ws.Range("K16").Value = -6.288343558282
$ws.Range("L16").Value = 53.132832080200
$ws.Range("M16").Value = -36.354166666666
$ws.Range("N16").Value = -87.443485408960

# Row 17
$ws.Range("C17").Value = 42
$ws.Range("D17").Value = 57
$ws.Range("E17").Value = -26.315789473684
$ws.Range("F17").Value = 225
$ws.Range("G17").Value = 263
$ws.Range("H17").Value = -14.448669201520
$ws.Range("I17").Value = 1150
$ws.Range("J17").Value = 1109
$ws.Range("K17").Value = 3.697024346257
$ws.Range("L17").Value = 23.922413793103
$ws.Range("M17").Value = 54.362416107382
$ws.Range("N17").Value = -46.336910872608

# Row 18
$ws.Range("C18").Value = 24
$ws.Range("D18").Value = 32
$ws.Range("E18").Value = -25
$ws.Range("F18").Value = 108
$ws.Range("G18").Value = 134
$ws.Range("H18").Value = -19.402985074626
$ws.Range("I18").Value = 623
$ws.Range("J18").Value = 689
$ws.Range("K18").Value = -9.579100145137
$ws.Range("L18").Value = 20.736434108527
$ws.Range("M18").Value = -45.778938207136
$ws.Range("N18").Value = -90.524714828897

# Row 19
$ws.Range("C19").Value = 138
$ws.Range("D19").Value = 123
$ws.Range("E19").Value = 12.195121951219
$ws.Range("F19").Value = 496
$ws.Range("G19").Value = 507
$ws.Range("H19").Value = -2.169625246548
$ws.Range("I19").Value = 2273
$ws.Range("J19").Value = 2413
$ws.Range("K19").Value = -5.801906340654
$ws.Range("L19").Value = 50.430178689609
$ws.Range("M19").Value = 30.332568807339
$ws.Range("N19").Value = -22.024013722126

# Row 20
$ws.Range("C20").Value = 27
$ws.Range("D20").Value = 24
$ws.Range("E20").Value = 12.5
$ws.Range("F20").Value = 131
$ws.Range("G20").Value = 121
$ws.Range("H20").Value = 8.264462809917
$ws.Range("I20").Value = 610
$ws.Range("J20").Value = 580
$ws.Range("K20").Value = 5.172413793103
$ws.Range("L20").Value = 61.375661375661
$ws.Range("M20").Value = -12.230215827338
$ws.Range("N20").Value = -92.810842663523

# Row 21
$ws.Range("C21").Value = 275
$ws.Range("D21").Value = 273
$ws.Range("E21").Value = 0.732600732600
$ws.Range("F21").Value = 1116
$ws.Range("G21").Value = 1158
$ws.Range("H21").Value = -3.626943005181
$ws.Range("I21").Value = 5366
$ws.Range("J21").Value = 5540
$ws.Range("K21").Value = -3.140794223826
$ws.Range("L21").Value = 40.544787847040
$ws.Range("M21").Value = -0.055876327062
$ws.Range("N21").Value = -78.743463793376

# Row 22
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 5
$ws.Range("E22").Value = -80
$ws.Range("F22").Value = 9
$ws.Range("G22").Value = 22
$ws.Range("H22").Value = -59.090909090909
$ws.Range("I22").Value = 57
$ws.Range("J22").Value = 67
$ws.Range("K22").Value = -14.925373134328
$ws.Range("L22").Value = 39.024390243902
$ws.Range("M22").Value = -39.361702127659

# Row 23
$ws.Range("C23").Value = 9
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = 80
$ws.Range("F23").Value = 35
$ws.Range("G23").Value = 37
$ws.Range("H23").Value = -5.405405405405
$ws.Range("I23").Value = 168
$ws.Range("J23").Value = 158
$ws.Range("K23").Value = 6.329113924050
$ws.Range("L23").Value = 3.067484662576
$ws.Range("M23").Value = 88.764044943820

# Row 24
$ws.Range("C24").Value = 302
$ws.Range("D24").Value = 273
$ws.Range("E24").Value = 10.622710622710
$ws.Range("F24").Value = 1312
$ws.Range("G24").Value = 1165
$ws.Range("H24").Value = 12.618025751073
$ws.Range("I24").Value = 5622
$ws.Range("J24").Value = 5122
$ws.Range("K24").Value = 9.761811792268
$ws.Range("L24").Value = 48.927152317880
$ws.Range("M24").Value = 41.826437941473

# Row 25
$ws.Range("C25").Value = 137
$ws.Range("D25").Value = 105
$ws.Range("E25").Value = 30.476190476190
$ws.Range("F25").Value = 470
$ws.Range("G25").Value = 393
$ws.Range("H25").Value = 19.592875318066
$ws.Range("I25").Value = 1938
$ws.Range("J25").Value = 1776
$ws.Range("K25").Value = 9.121621621621
$ws.Range("L25").Value = 27.332457293035
$ws.Range("M25").Value = -10.360777058279

# Row 26
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = -50
$ws.Range("F26").Value = 16
$ws.Range("G26").Value = 23
$ws.Range("H26").Value = -30.434782608695
$ws.Range("I26").Value = 105
$ws.Range("J26").Value = 119
$ws.Range("K26").Value = -11.764705882352
$ws.Range("L26").Value = -7.079646017699

# Row 27
$ws.Range("C27").Value = 26
$ws.Range("D27").Value = 17
$ws.Range("E27").Value = 52.941176470588
$ws.Range("F27").Value = 62
$ws.Range("G27").Value = 52
$ws.Range("H27").Value = 19.230769230769
$ws.Range("I27").Value = 205
$ws.Range("J27").Value = 216
$ws.Range("K27").Value = -5.092592592592
$ws.Range("L27").Value = 14.525139664804

# Row 28
$ws.Range("C28").Value = 3
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 200
$ws.Range("F28").Value = 14
$ws.Range("G28").Value = 20
$ws.Range("H28").Value = -30
$ws.Range("I28").Value = 57
$ws.Range("J28").Value = 73
$ws.Range("K28").Value = -21.917808219178
$ws.Range("L28").Value = -19.718309859154
$ws.Range("M28").Value = -10.9375
$ws.Range("N28").Value = -77.108433734939

# Row 29
$ws.Range("C29").Value = 3
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = 200
$ws.Range("F29").Value = 13
$ws.Range("G29").Value = 11
$ws.Range("H29").Value = 18.181818181818
$ws.Range("I29").Value = 48
$ws.Range("J29").Value = 53
$ws.Range("K29").Value = -9.433962264150
$ws.Range("L29").Value = -25
$ws.Range("M29").Value = -9.433962264150
$ws.Range("N29").Value = -78.475336322869

# Row 30
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = 6
$ws.Range("E30").Value = -83.333333333333
$ws.Range("F30").Value = 9
$ws.Range("G30").Value = 14
$ws.Range("H30").Value = -35.714285714285
$ws.Range("I30").Value = 21
$ws.Range("J30").Value = 48
$ws.Range("K30").Value = -56.25
$ws.Range("L30").Value = 23.529411764705
